# Add a new entry to the "daily problem" tracking sheet for LeetCode 887
# (hard dp), marked as done, dated 2019-03-29 (Excel serial 43553).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the date formatting from the previous date cell (A25) so the new
# date cell (A26) picks up the existing date number format instead of
# creating a brand-new style.
$ws.Range("A25").Copy()
$ws.Range("A26").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = 0

# Fill in the new row's data.
$ws.Range("A26").Value = 43553
$ws.Range("B26").Value = "887 dp"
$ws.Range("F26").Value = "done"

# Match the author's final selection (cell below the newly entered row).
$ws.Range("B27").Select()
